$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Underline the "Type" header cell (B1)
$ws.Range("B1").Font.Underline = $true

# Correct the data type listed for example_spike_count (row 18) and
# example_segment (row 19): both are actually "cell" arrays, not "double"
$ws.Range("B18").Value = "cell"
$ws.Range("B19").Value = "cell"

# Move the active selection back to B1
$ws.Range("B1").Select() | Out-Null
